$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Insert the two new rows first (so row numbers shift naturally) ---
# Row 2: new "Label" field (Welcome banner)
$ws1.Rows.Item(2).Insert()
# Row 4 (after insert above, "sex" question is now row 5, so GPS question goes at row 4,
# right after the Patient-name question which will be row 3)
$ws1.Rows.Item(4).Insert()

# --- Row 2: Label / Welcome ---
$ws1.Range("A2").Value = "Welcome to the outbreak questionnaire"
$ws1.Range("B2").Value = "Welcome"
$ws1.Range("D2").Value = "label1"
$ws1.Range("E2").Value = "Label"
$ws1.Range("C2:J2").Clear()
$ws1.Range("A2:D2").Font.Bold = $false

# --- Row 3: existing patient-name question, now gets a Description ---
$ws1.Range("C3").Value = "Use this space to provide additional instructions to the interviewer"

# --- Row 4: new GPS / Interview GPS question ---
$ws1.Range("A4").Value = "Where is the interview occurring?"
$ws1.Range("B4").Value = "Interview GPS"
$ws1.Range("D4").Value = "interview_gps"
$ws1.Range("E4").Value = "GPS"
$ws1.Range("F4").Value = $false
$ws1.Range("C4:J4").Clear()
$ws1.Range("F4").Value = $false

# --- Remove the old "Please enter ..." placeholder descriptions that are no longer used ---
$ws1.Range("C5").ClearContents()   # age
$ws1.Range("C6").ClearContents()   # sex
$ws1.Range("C7").ClearContents()   # pregnant
$ws1.Range("C8").ClearContents()   # onset date
$ws1.Range("C9").ClearContents()   # onset time
$ws1.Range("C10").ClearContents()  # eaten foods

# --- selection ---
$ws1.Range("E4").Select() | Out-Null

# --- DataTypes sheet: append GPS and Label as new selectable types ---
$ws5 = $wb.Worksheets.Item("DataTypes")
$ws5.Range("A9").Value = "GPS"
$ws5.Range("A10").Value = "Label"
$ws5.Range("A9").Select() | Out-Null
